# Auto-generated: update computed pl_mw results for "case with 380 kV done"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B2" = 0.674334574529496
    "D2" = 0.165699909200427
    "E2" = 0.1709213590108387
    "F2" = 1.652445559581004
    "G2" = 0.002484681466328701
    "I2" = 1.221848163755823
    "J2" = 0.2209342131489471
    "K2" = 0.4599897358796738
    "L2" = 0.2430947639491876
    "O2" = 4.112535607678694
    "B3" = 0.6465430905678318
    "D3" = 0.1652084267701852
    "E3" = 0.1719248794348829
    "F3" = 1.661099994329319
    "G3" = 0.002487215116926315
    "I3" = 1.234682631242698
    "J3" = 0.2224605919060992
    "K3" = 0.4052526514537078
    "L3" = 0.2319176901096398
    "O3" = 4.138691301653296
    "B4" = 0.6296693325681133
    "D4" = 0.1649548143632131
    "E4" = 0.1725870779272398
    "F4" = 1.667209193760712
    "G4" = 0.002488854977137414
    "I4" = 1.243081333061472
    "J4" = 0.2234508350012261
    "K4" = 0.3715842841881454
    "L4" = 0.2251290745373211
    "O4" = 4.156845794165704
    "B5" = 0.6228416270937771
    "D5" = 0.1648636280581854
    "E5" = 0.1728685289161067
    "F5" = 1.66989888730707
    "G5" = 0.002489544465576963
    "I5" = 1.246634257930797
    "J5" = 0.2238677287698897
    "K5" = 0.3578500094317008
    "L5" = 0.2223815064890431
    "O5" = 4.164770822221215
    "B6" = 0.6217108388959787
    "D6" = 0.1648492227955884
    "E6" = 0.1729159649500378
    "F6" = 1.670357601560802
    "G6" = 0.002489660238651603
    "I6" = 1.247232095086893
    "J6" = 0.2239377613870859
    "K6" = 0.3555686142129559
    "L6" = 0.2219264192179509
    "O6" = 4.166118597159297
    "B7" = 0.6295770545827111
    "D7" = 0.1649535352752025
    "E7" = 0.1725908266742104
    "F7" = 1.667244657330059
    "G7" = 0.002488864189827436
    "I7" = 1.243128721044965
    "J7" = 0.2234564032429915
    "K7" = 0.3713991151004734
    "L7" = 0.2250919432912184
    "O7" = 4.156950540141167
    "B8" = 0.6647129359516555
    "D8" = 0.1655204816203479
    "E8" = 0.1712578333889603
    "F8" = 1.655264624112718
    "G8" = 0.002485537633317636
    "I8" = 1.226165953683694
    "J8" = 0.2214495146920435
    "K8" = 0.4411292401389346
    "L8" = 0.2392256354014393
    "O8" = 4.121119417918663
    "B9" = 0.7351006669343576
    "D9" = 0.1670122462673902
    "E9" = 0.1690080009632551
    "F9" = 1.638076721055718
    "G9" = 0.002479679427796303
    "I9" = 1.197011511540214
    "J9" = 0.2179337389874222
    "K9" = 0.5773668417247109
    "L9" = 0.26752282070197
    "O9" = 4.067470156975929
    "B10" = 0.7876943993611007
    "D10" = 0.1683372668864109
    "E10" = 0.1675755677439579
    "F10" = 1.629285867523151
    "G10" = 0.002475776932463228
    "I10" = 1.178092653380386
    "J10" = 0.2156049738108184
    "K10" = 0.6771230756575335
    "L10" = 0.2886594368673769
    "O10" = 4.038178251005149
    "B11" = 0.8118065059259152
    "D11" = 0.1689892433045728
    "E11" = 0.1669714859446856
    "F11" = 1.626118576672212
    "G11" = 0.002474087929171029
    "I11" = 1.170028095995338
    "J11" = 0.2146004301255294
    "K11" = 0.7224255489542486
    "L11" = 0.2983488494975717
    "O11" = 4.02705014297581
    "B12" = 0.8209634595758075
    "D12" = 0.1692431576844484
    "E12" = 0.1667495473962202
    "F12" = 1.625038681879957
    "G12" = 0.002473460688366184
    "I12" = 1.167052101816378
    "J12" = 0.2142278941525569
    "K12" = 0.7395686097115401
    "L12" = 0.3020284839329292
    "O12" = 4.023152034867252
    "B13" = 0.8189901900403527
    "D13" = 0.1691881610715882
    "E13" = 0.16679704312285
    "F13" = 1.625265943969836
    "G13" = 0.002473595227629998
    "I13" = 1.16768957206402
    "J13" = 0.2143077770607785
    "K13" = 0.7358770931320748
    "L13" = 0.3012355468322596
    "O13" = 4.023977514835792
    "B14" = 0.8125593322636178
    "D14" = 0.1690099925009818
    "E14" = 0.1669530904870982
    "F14" = 1.62602733896027
    "G14" = 0.002474036078727328
    "I14" = 1.169781698476761
    "J14" = 0.2145696238998303
    "K14" = 0.7238361646349745
    "L14" = 0.2986513667461281
    "O14" = 4.026723113166611
    "B15" = 0.8086236414938242
    "D15" = 0.1689017724239577
    "E15" = 0.1670495607524352
    "F15" = 1.626509273173369
    "G15" = 0.002474307719058766
    "I15" = 1.171073327802748
    "J15" = 0.2147310360247502
    "K15" = 0.7164591581882291
    "L15" = 0.29706983821508
    "O15" = 4.028446004834393
    "B16" = 0.7861223218820044
    "D16" = 0.1682956445995742
    "E16" = 0.1676160006415159
    "F16" = 1.629509584242086
    "G16" = 0.002475889043546502
    "I16" = 1.178630590293228
    "J16" = 0.2156717245227355
    "K16" = 0.6741608223811113
    "L16" = 0.2880276869332903
    "O16" = 4.038949702529607
    "B17" = 0.7723659031189811
    "D17" = 0.1679363733528731
    "E17" = 0.1679756534511068
    "F17" = 1.631563121508961
    "G17" = 0.002476881188335892
    "I17" = 1.183405460137209
    "J17" = 0.2162628332915932
    "K17" = 0.6481917577816034
    "L17" = 0.2824994871156719
    "O17" = 4.045956030964817
    "B18" = 0.7644712112860361
    "D18" = 0.1677343657809729
    "E18" = 0.1681869921631964
    "F18" = 1.63282255313657
    "G18" = 0.002477459966489043
    "I18" = 1.186202830027504
    "J18" = 0.2166079848688445
    "K18" = 0.6332478323099053
    "L18" = 0.2793268152774289
    "O18" = 4.050192674009168
    "B19" = 0.7618012535417336
    "D19" = 0.16766676728097
    "E19" = 0.1682593172354174
    "F19" = 1.633262424990768
    "G19" = 0.002477657327813496
    "I19" = 1.187158729811422
    "J19" = 0.2167257343760866
    "K19" = 0.6281868660080931
    "L19" = 0.2782538123276623
    "O19" = 4.051662646155989
    "B20" = 0.7738284768787196
    "D20" = 0.1679741389891305
    "E20" = 0.1679369047134216
    "F20" = 1.631336416897177
    "G20" = 0.002476774732785742
    "I20" = 1.182891889860514
    "J20" = 0.2161993747134083
    "K20" = 0.6509569614781299
    "L20" = 0.2830872502474477
    "O20" = 4.045188793893232
    "B21" = 0.8144475246783145
    "D21" = 0.1690621346879624
    "E21" = 0.1669070708350961
    "F21" = 1.625800456761141
    "G21" = 0.002473906255340727
    "I21" = 1.16916507679522
    "J21" = 0.2144924999016508
    "K21" = 0.7273732101999144
    "L21" = 0.2994101207906965
    "O21" = 4.025908092308754
    "B22" = 0.8411469580162532
    "D22" = 0.1698141259796202
    "E22" = 0.1662737226290449
    "F22" = 1.622878811849404
    "G22" = 0.002472103486581343
    "I22" = 1.160647738035735
    "J22" = 0.2134227802115186
    "K22" = 0.7772452651567505
    "L22" = 0.3101389542432429
    "O22" = 4.015148117279097
    "B23" = 0.8268832288095496
    "D23" = 0.1694090468760407
    "E23" = 0.1666081264581543
    "F23" = 1.624374462169484
    "G23" = 0.002473059093555541
    "I23" = 1.165152074476023
    "J23" = 0.2139895238193079
    "K23" = 0.7506343506663598
    "L23" = 0.3044072728207539
    "O23" = 4.020722476695255
    "B24" = 0.773167203448196
    "D24" = 0.1679570509935147
    "E24" = 0.1679544087828262
    "F24" = 1.631438664506746
    "G24" = 0.00247682283521011
    "I24" = 1.183123912276297
    "J24" = 0.2162280477785754
    "K24" = 0.6497068563168966
    "L24" = 0.2828215051893181
    "O24" = 4.045535011954399
    "B25" = 0.7159025181755965
    "D25" = 0.1665682766745817
    "E25" = 0.1695778079086612
    "F25" = 1.642052195802528
    "G25" = 0.002481193429915567
    "I25" = 1.204459062131409
    "J25" = 0.2188400810185084
    "K25" = 0.5773668417247109
    "L25" = 0.2598062588371732
    "O25" = 4.080205421433988
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
